# Auto-generated edit script: updates Goblin Profits market-data snapshot values
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR crafting-profit worksheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 93.63636
$ws.Range("I8").Value = 93.63636
$ws.Range("K8").Value = 280.90908
$ws.Range("M8").Value = -141.90908
$ws.Range("H11").Value = 72.77778000000001
$ws.Range("I11").Value = 72.77778000000001
$ws.Range("K11").Value = 72.77778000000001
$ws.Range("M11").Value = 67.22221999999999
$ws.Range("H19").Value = 602.5
$ws.Range("I19").Value = 744.2
$ws.Range("K19").Value = 744.2
$ws.Range("M19").Value = -569.2
$ws.Range("H28").Value = 6623.706
$ws.Range("I28").Value = 8245.923000000001
$ws.Range("J28").Value = 1351.5
$ws.Range("K28").Value = 8245.923000000001
$ws.Range("L28").Value = 1351.5
$ws.Range("M28").Value = -7760.923000000001
$ws.Range("N28").Value = -2321.5
$ws.Range("H41").Value = 986.61536
$ws.Range("I41").Value = 837
$ws.Range("J41").Value = 1114.8572
$ws.Range("K41").Value = 837
$ws.Range("L41").Value = 1114.8572
$ws.Range("M41").Value = -397
$ws.Range("N41").Value = -1994.8572
$ws.Range("H106").Value = 4253.8335
$ws.Range("I106").Value = 4253.8335
$ws.Range("K106").Value = 4253.8335
$ws.Range("M106").Value = -3622.8335
$ws.Range("H129").Value = 1876.2222
$ws.Range("J129").Value = 2999.2
$ws.Range("L129").Value = 8997.599999999999
$ws.Range("N129").Value = -18997.6
$ws.Range("H138").Value = 3645.1316
$ws.Range("I138").Value = 1912.9333
$ws.Range("K138").Value = 5738.7999
$ws.Range("M138").Value = -598.7999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 175
$ws.Range("I4").Value = 175
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 175
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -59
$ws.Range("N4").Value = $null
$ws.Range("H8").Value = 5135001.5
$ws.Range("I8").Value = 6836666.5
$ws.Range("K8").Value = 6836666.5
$ws.Range("M8").Value = -6836522.5
$ws.Range("H50").Value = 8574
$ws.Range("I50").Value = 839.6667
$ws.Range("J50").Value = 12441.167
$ws.Range("K50").Value = 839.6667
$ws.Range("L50").Value = 12441.167
$ws.Range("M50").Value = -125.6667
$ws.Range("N50").Value = -13869.167
$ws.Range("H61").Value = 3030.2341
$ws.Range("I61").Value = 3022.8096
$ws.Range("K61").Value = 3022.8096
$ws.Range("M61").Value = -2810.8096
$ws.Range("H97").Value = 262.4
$ws.Range("I97").Value = 254.57895
$ws.Range("K97").Value = 254.57895
$ws.Range("M97").Value = 241.42105
$ws.Range("H122").Value = 11115137
$ws.Range("J122").Value = 4398.6665
$ws.Range("L122").Value = 13195.9995
$ws.Range("N122").Value = -18095.9995
$ws.Range("H136").Value = 3030.2341
$ws.Range("I136").Value = 3022.8096
$ws.Range("K136").Value = 9068.4288
$ws.Range("M136").Value = -6518.4288

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 4120.8335
$ws.Range("I99").Value = 2333.3333
$ws.Range("J99").Value = 4716.6665
$ws.Range("K99").Value = 2333.3333
$ws.Range("L99").Value = 4716.6665
$ws.Range("M99").Value = -835.3332999999998
$ws.Range("N99").Value = -7712.6665
$ws.Range("H134").Value = 2098.7673
$ws.Range("I134").Value = 2048.6667
$ws.Range("K134").Value = 6146.000100000001
$ws.Range("M134").Value = -3611.000100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2220.3
$ws.Range("I16").Value = 1911.4445
$ws.Range("K16").Value = 1911.4445
$ws.Range("M16").Value = -1624.4445
$ws.Range("H31").Value = 4189.6
$ws.Range("I31").Value = 2562.4
$ws.Range("K31").Value = 2562.4
$ws.Range("M31").Value = -2267.4
$ws.Range("H34").Value = 4189.6
$ws.Range("I34").Value = 2562.4
$ws.Range("K34").Value = 2562.4
$ws.Range("M34").Value = -2360.4
$ws.Range("H58").Value = 1387.2354
$ws.Range("I58").Value = 1421.8572
$ws.Range("J58").Value = 1225.6666
$ws.Range("K58").Value = 1421.8572
$ws.Range("L58").Value = 1225.6666
$ws.Range("M58").Value = -1218.8572
$ws.Range("N58").Value = -1631.6666
$ws.Range("H97").Value = 123089.5
$ws.Range("J97").Value = 185000
$ws.Range("L97").Value = 185000
$ws.Range("N97").Value = -186982
$ws.Range("H113").Value = 2220.3
$ws.Range("I113").Value = 1911.4445
$ws.Range("K113").Value = 1911.4445
$ws.Range("M113").Value = 258.5554999999999
$ws.Range("H136").Value = 1387.2354
$ws.Range("I136").Value = 1421.8572
$ws.Range("J136").Value = 1225.6666
$ws.Range("K136").Value = 4265.571599999999
$ws.Range("L136").Value = 3676.9998
$ws.Range("M136").Value = -1715.571599999999
$ws.Range("N136").Value = -8776.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2117.1875
$ws.Range("I5").Value = 617
$ws.Range("J5").Value = 4617.5
$ws.Range("K5").Value = 1851
$ws.Range("L5").Value = 13852.5
$ws.Range("M5").Value = -1739
$ws.Range("N5").Value = -14076.5
$ws.Range("H46").Value = 44848484
$ws.Range("I46").Value = 69696970
$ws.Range("K46").Value = 209090910
$ws.Range("M46").Value = -209090819
$ws.Range("H121").Value = 720.75
$ws.Range("J121").Value = 806.5
$ws.Range("L121").Value = 2419.5
$ws.Range("N121").Value = -5039.5
$ws.Range("H132").Value = 2499.8
$ws.Range("I132").Value = 2166.3333
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 19496.9997
$ws.Range("L132").Value = 27000
$ws.Range("M132").Value = -16966.9997
$ws.Range("N132").Value = -32060
$ws.Range("H133").Value = 5165.2
$ws.Range("I133").Value = 3407.0908
$ws.Range("K133").Value = 10221.2724
$ws.Range("M133").Value = -5161.2724
$ws.Range("H134").Value = 2041.8
$ws.Range("I134").Value = 1432.8462
$ws.Range("K134").Value = 4298.5386
$ws.Range("M134").Value = 771.4614000000001
$ws.Range("H135").Value = 2117.1875
$ws.Range("I135").Value = 617
$ws.Range("J135").Value = 4617.5
$ws.Range("K135").Value = 5553
$ws.Range("L135").Value = 41557.5
$ws.Range("M135").Value = -3018
$ws.Range("N135").Value = -46627.5
$ws.Range("H136").Value = 4499.6665
$ws.Range("I136").Value = 4499.6665
$ws.Range("K136").Value = 13498.9995
$ws.Range("M136").Value = -8398.999500000002
$ws.Range("H137").Value = 3820.3333
$ws.Range("I137").Value = 3820.3333
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 11460.9999
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -6360.999899999999
$ws.Range("N137").Value = $null
$ws.Range("H138").Value = 3543.0715
$ws.Range("I138").Value = 2560.3
$ws.Range("J138").Value = 6000
$ws.Range("K138").Value = 7680.900000000001
$ws.Range("L138").Value = 18000
$ws.Range("M138").Value = -2540.900000000001
$ws.Range("N138").Value = -28280
$ws.Range("H140").Value = 3693.2144
$ws.Range("I140").Value = 3693.2144
$ws.Range("K140").Value = 11079.6432
$ws.Range("M140").Value = -5899.643199999999
$ws.Range("H141").Value = 1165.25
$ws.Range("I141").Value = 1165.25
$ws.Range("K141").Value = 3495.75
$ws.Range("M141").Value = 1684.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 1256750.8
$ws.Range("I10").Value = 2503501.5
$ws.Range("J10").Value = 10000
$ws.Range("K10").Value = 2503501.5
$ws.Range("L10").Value = 10000
$ws.Range("M10").Value = -2503332.5
$ws.Range("N10").Value = -10338
$ws.Range("H113").Value = 47627920
$ws.Range("I113").Value = 333334660
$ws.Range("K113").Value = 333334660
$ws.Range("M113").Value = -333332490

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1417.7273
$ws.Range("I22").Value = 2608.9
$ws.Range("J22").Value = 1067.3823
$ws.Range("K22").Value = 2608.9
$ws.Range("L22").Value = 1067.3823
$ws.Range("M22").Value = -2313.9
$ws.Range("N22").Value = -1657.3823
$ws.Range("H27").Value = 1417.7273
$ws.Range("I27").Value = 2608.9
$ws.Range("J27").Value = 1067.3823
$ws.Range("K27").Value = 2608.9
$ws.Range("L27").Value = 1067.3823
$ws.Range("M27").Value = -2501.9
$ws.Range("N27").Value = -1281.3823
$ws.Range("H100").Value = 9400.200000000001
$ws.Range("I100").Value = 6999.5
$ws.Range("K100").Value = 6999.5
$ws.Range("M100").Value = -6458.5
$ws.Range("H122").Value = 5931.56
$ws.Range("I122").Value = 5480.2383
$ws.Range("K122").Value = 16440.7149
$ws.Range("M122").Value = -13990.7149
$ws.Range("H132").Value = 5530.6924
$ws.Range("I132").Value = 5987.375
$ws.Range("J132").Value = 4800
$ws.Range("K132").Value = 17962.125
$ws.Range("L132").Value = 14400
$ws.Range("M132").Value = -15432.125
$ws.Range("N132").Value = -19460

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 6413.857
$ws.Range("I96").Value = 7849.25
$ws.Range("K96").Value = 7849.25
$ws.Range("M96").Value = -6476.25
$ws.Range("H100").Value = 2351.5
$ws.Range("I100").Value = 700
$ws.Range("J100").Value = 4003
$ws.Range("K100").Value = 1400
$ws.Range("L100").Value = 8006
$ws.Range("M100").Value = -859
$ws.Range("N100").Value = -9088
$ws.Range("H107").Value = 4493.4
$ws.Range("I107").Value = 3848
$ws.Range("J107").Value = 5999.3335
$ws.Range("K107").Value = 11544
$ws.Range("L107").Value = 17998.0005
$ws.Range("M107").Value = -9624
$ws.Range("N107").Value = -21838.0005
$ws.Range("H113").Value = 930.3946999999999
$ws.Range("I113").Value = 825.7586
$ws.Range("J113").Value = 1267.5555
$ws.Range("K113").Value = 2477.2758
$ws.Range("L113").Value = 3802.6665
$ws.Range("M113").Value = -307.2757999999999
$ws.Range("N113").Value = -8142.666499999999
